$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.053.89"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.870.45"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.25%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5108"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3884"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.82%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08351"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.114"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.39%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.208"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.882.61"
$ws.Range("E12").Value = "  +1.06%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.59"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.230"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.00001098"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "90.59"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06667"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "17.70"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.993"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "28.108.46"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +0.39%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.250"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "159.03"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.463"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.60%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.53"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.54%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "124.70"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.1055"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.030"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.809"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.37%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.593"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.22%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "9.534"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.02445"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.06532"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.2184"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.72%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.192"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.06%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.6461"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.964"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.78%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.219"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("E41").Value = "  +0.66%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.6101"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.30%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "13.06"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.280"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.663"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.32%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.001"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +1.33%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "120.34"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06877"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "77.76"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.35%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.1386"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.96%  "
